{"js": "// The attendance table is the second table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[1];\n\n// --- 1. Merge the split \"08\" / \"/03/2023\" runs in the 08/03/2023 10:30 row\n//        (row index 6, 0-based) into a single run, same end text, same\n//        paragraph formatting.\nconst dateCell = table.getCell(6, 0);\ndateCell.body.load(\"paragraphs\");\nawait context.sync();\n\nconst dateParagraph = dateCell.body.paragraphs.items[0];\ndateParagraph.getRange().insertText(\"08/03/2023\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2. Append the two new attendance rows at the bottom of the table.\ntable.addRows(Word.InsertLocation.end, 2, [\n  [\"09/03/2023\", \"14:30\", \"1.5\", \"Observer\", \"Driver\"],\n  [\"09/03/2023\", \"15:00\", \"1.5\", \"Driver\", \"Observer\"]\n]);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The attendance table is the second table in the document.\n$tbl = $d.Tables.Item(2)\n\n# --- 1. Merge the split \"08\" / \"/03/2023\" runs in the 08/03/2023 10:30 row\n#        into a single run of text. A plain Range.Text assignment only\n#        overwrites the first run when a cell holds more than one run, so\n#        use Find/Replace (scoped to that cell) to rewrite the whole cell\n#        content atomically, exactly like retyping the date in Word would.\n$dateCell = $tbl.Cell(7, 1)\n$findRange = $dateCell.Range\n$findRange.Find.ClearFormatting()\n$findRange.Find.Execute(\"08/03/2023\", $false, $false, $false, $false, $false, $true, 1, $false, \"08/03/2023\", 2) | Out-Null\n\n# --- 2. Append the two new attendance rows at the bottom of the table.\n$tbl.Rows.Add() | Out-Null\n$tbl.Rows.Add() | Out-Null\n\n$row1 = $tbl.Rows.Count - 1\n$row2 = $tbl.Rows.Count\n\n$tbl.Cell($row1, 1).Range.Text = \"09/03/2023\"\n$tbl.Cell($row1, 2).Range.Text = \"14:30\"\n$tbl.Cell($row1, 3).Range.Text = \"1.5\"\n$tbl.Cell($row1, 4).Range.Text = \"Observer\"\n$tbl.Cell($row1, 5).Range.Text = \"Driver\"\n\n$tbl.Cell($row2, 1).Range.Text = \"09/03/2023\"\n$tbl.Cell($row2, 2).Range.Text = \"15:00\"\n$tbl.Cell($row2, 3).Range.Text = \"1.5\"\n$tbl.Cell($row2, 4).Range.Text = \"Driver\"\n$tbl.Cell($row2, 5).Range.Text = \"Observer\"\n"}
